$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 531, pushing existing rows 531-619
# down to 533-621.
$ws.Rows("531:532").Insert()

# --- Populate new row 531 ---
$ws.Cells.Item(531, 1).Value = 11
$ws.Cells.Item(531, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(531, 3).Value = "Bíobío"
$ws.Cells.Item(531, 4).Value = 44951
$ws.Cells.Item(531, 5).Value = 8
$ws.Cells.Item(531, 6).Value = 100112020
$ws.Cells.Item(531, 7).Value = "Tomate"
$ws.Cells.Item(531, 8).Value = "Semiduro"
$ws.Cells.Item(531, 9).Value = "Primera"
$ws.Cells.Item(531, 10).Value = 500
$ws.Cells.Item(531, 11).Value = 8000
$ws.Cells.Item(531, 12).Value = 9000
$ws.Cells.Item(531, 13).Value = 8600
$ws.Cells.Item(531, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(531, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(531, 16).Value = 478
$ws.Cells.Item(531, 17).Value = 18
$ws.Cells.Item(531, 18).Value = "Hortaliza"

# --- Populate new row 532 ---
$ws.Cells.Item(532, 1).Value = 11
$ws.Cells.Item(532, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(532, 3).Value = "Bíobío"
$ws.Cells.Item(532, 4).Value = 44951
$ws.Cells.Item(532, 5).Value = 8
$ws.Cells.Item(532, 6).Value = 100112020
$ws.Cells.Item(532, 7).Value = "Tomate"
$ws.Cells.Item(532, 8).Value = "Semiduro"
$ws.Cells.Item(532, 9).Value = "Segunda"
$ws.Cells.Item(532, 10).Value = 300
$ws.Cells.Item(532, 11).Value = 7000
$ws.Cells.Item(532, 12).Value = 7000
$ws.Cells.Item(532, 13).Value = 7000
$ws.Cells.Item(532, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(532, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(532, 16).Value = 389
$ws.Cells.Item(532, 17).Value = 18
$ws.Cells.Item(532, 18).Value = "Hortaliza"
